$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-66 down to 55-67.
$ws.Rows.Item(54).Insert()

# Populate the new row 54 with the new weekly record.
$ws.Range("A54").Value = 7
$ws.Range("B54").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C54").Value = "Ñuble"
$ws.Range("D54").Value = 44588
$ws.Range("D54").NumberFormat = $ws.Range("D55").NumberFormat
$ws.Range("E54").Value = 16
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100108
$ws.Range("H54").Value = "Tropicales y subtropicales"
$ws.Range("I54").Value = 100108002
$ws.Range("J54").Value = "Mango"
$ws.Range("K54").Value = "Sin especificar"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 60
$ws.Range("N54").Value = 7000
$ws.Range("O54").Value = 7500
$ws.Range("P54").Value = 7250
$ws.Range("Q54").Value = "$/bandeja 4 kilos"
$ws.Range("R54").Value = "Perú"
$ws.Range("S54").Value = 1812
$ws.Range("T54").Value = 4
